# Apply revised figures to Table 2 and Table 3 per commit
# "Revised figures, tables, manuscript and supplement".

$wb = $excel.ActiveWorkbook

# --- Table 2: "Mountain sport type" breakdown (B3) ---
$wsTable2 = $wb.Worksheets.Item("Table 2")
$mountainSportType = "alpine skiing/snowboarding: 59% (n = 180)`nski touring/freeride: 2.9% (n = 9)`ncross-country skiing: 5.5% (n = 17)`nsledding: 3.9% (n = 12)`nice climbing: 0.33% (n = 1)`nhiking: 5.5% (n = 17)`nclimbing: 3.6% (n = 11)`nmountaineering: 0.65% (n = 2)`nbiking: 16% (n = 48)`nair sport: 0.33% (n = 1)`nwater sport: 1.6% (n = 5)`nother: 1.3% (n = 4)`nn = 307"
$wsTable2.Range("B3").Value = $mountainSportType

# --- Table 3: anxiety / depression figures ---
$wsTable3 = $wb.Worksheets.Item("Table 3")
$wsTable3.Range("B4").Value = "2.9% (n = 9)"
$wsTable3.Range("A6").Value = "Clinically relevant depression symptoms (PHQ-9 ≥10)"
$wsTable3.Range("B6").Value = "7.2% (n = 22)"
